# Fruta / hortaliza, semanal
# Inserts 3 new weekly records at the top of the "Ají" data block (rows 78-80),
# pushing the existing records (old rows 78-143) down to rows 81-146.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 78..143 down by 3 rows (inserting 3 blank rows at row 78).
$ws.Rows.Item(78).Resize(3).Insert()

# New row 78
$ws.Cells.Item(78, 1).Value = 8
$ws.Cells.Item(78, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(78, 3).Value = "Coquimbo"
$ws.Cells.Item(78, 4).Value = 44539
$ws.Cells.Item(78, 5).Value = 4
$ws.Cells.Item(78, 6).Value = 100112021
$ws.Cells.Item(78, 7).Value = "Ají"
$ws.Cells.Item(78, 8).Value = "Inferno"
$ws.Cells.Item(78, 9).Value = "Primera"
$ws.Cells.Item(78, 10).Value = 500
$ws.Cells.Item(78, 11).Value = 15000
$ws.Cells.Item(78, 12).Value = 16000
$ws.Cells.Item(78, 13).Value = 15500
$ws.Cells.Item(78, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(78, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(78, 16).Value = 1292
$ws.Cells.Item(78, 17).Value = 12
$ws.Cells.Item(78, 18).Value = "Hortaliza"

# New row 79
$ws.Cells.Item(79, 1).Value = 8
$ws.Cells.Item(79, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(79, 3).Value = "Coquimbo"
$ws.Cells.Item(79, 4).Value = 44539
$ws.Cells.Item(79, 5).Value = 4
$ws.Cells.Item(79, 6).Value = 100112021
$ws.Cells.Item(79, 7).Value = "Ají"
$ws.Cells.Item(79, 8).Value = "Inferno"
$ws.Cells.Item(79, 9).Value = "Primera"
$ws.Cells.Item(79, 10).Value = 500
$ws.Cells.Item(79, 11).Value = 18000
$ws.Cells.Item(79, 12).Value = 19000
$ws.Cells.Item(79, 13).Value = 18500
$ws.Cells.Item(79, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(79, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(79, 16).Value = 1233
$ws.Cells.Item(79, 17).Value = 15
$ws.Cells.Item(79, 18).Value = "Hortaliza"

# New row 80
$ws.Cells.Item(80, 1).Value = 8
$ws.Cells.Item(80, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(80, 3).Value = "Coquimbo"
$ws.Cells.Item(80, 4).Value = 44539
$ws.Cells.Item(80, 5).Value = 4
$ws.Cells.Item(80, 6).Value = 100112021
$ws.Cells.Item(80, 7).Value = "Ají"
$ws.Cells.Item(80, 8).Value = "Inferno"
$ws.Cells.Item(80, 9).Value = "Segunda"
$ws.Cells.Item(80, 10).Value = 300
$ws.Cells.Item(80, 11).Value = 8000
$ws.Cells.Item(80, 12).Value = 9000
$ws.Cells.Item(80, 13).Value = 8500
$ws.Cells.Item(80, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(80, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(80, 16).Value = 708
$ws.Cells.Item(80, 17).Value = 12
$ws.Cells.Item(80, 18).Value = "Hortaliza"
